$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 115; this shifts existing rows 115-150 down to 116-151
# and automatically extends the used range / dimension to R151.
$ws.Rows(115).Insert()

# Populate the newly inserted row 115 with the new weekly price record.
$ws.Range("A115").Value = 4
$ws.Range("B115").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C115").Value = "Los Lagos"
$ws.Range("D115").Value = 45093
$ws.Range("E115").Value = 10
$ws.Range("F115").Value = 100112031
$ws.Range("G115").Value = "Poroto verde"
$ws.Range("H115").Value = "Magnum"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 45
$ws.Range("K115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("M115").Value = 30000
$ws.Range("N115").Value = "$/malla 25 kilos"
$ws.Range("O115").Value = "Perú"
$ws.Range("P115").Value = 1200
$ws.Range("Q115").Value = 25
$ws.Range("R115").Value = "Hortaliza"
